$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (Beteckning) to determine the data extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 431 }

# Column C holds the "Förändrad" (changed) date, stored as serial 45186 (2023-09-17).
# Update every data row (2..lastRow) to serial 45188 (2023-09-19), preserving formatting.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
